# Generate Report for Handback
#
# This script updates the localization-status workbook to record the
# result of a handback run for af758db3-7343-497e-8cc8-91cff17baa62:
#  - the "zh-cn" sheet gets its Latest Target File / Latest Handback File /
#    Latest Handback DateTime / Error Detail columns filled in for row 6
#  - the "de-de" sheet gets the same treatment (with its own handback time)
#  - both sheets widen columns I, J and P to fit the new content and add
#    a hyperlink on the newly populated "Latest Target File" cell

$wb = $excel.ActiveWorkbook

$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ca6f03f575e62f26e3ca32f6798e551929d24152/e2e/af758db3-7343-497e-8cc8-91cff17baa62.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/241ab3f8dd3ac37ab8c934afd459b2653fca7e1b/e2e/af758db3-7343-497e-8cc8-91cff17baa62.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ca6f03f575e62f26e3ca32f6798e551929d24152/e2e/af758db3-7343-497e-8cc8-91cff17baa62.md."

function Update-LocaleSheet($sheetName, $handbackDateTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Widen the columns that now hold longer content.
    $ws.Columns.Item(9).ColumnWidth = 39.17
    $ws.Columns.Item(10).ColumnWidth = 39.17
    $ws.Columns.Item(16).ColumnWidth = 39.17

    # Latest Target File
    $ws.Range("I6").Value = "af758db3-7343-497e-8cc8-91cff17baa62.md"

    # Latest Handback File
    $ws.Range("J6").Value = "af758db3-7343-497e-8cc8-91cff17baa62.b449ec60f18b373a1e1f7a9216650947cd6f624a.$sheetName.xlf"

    # Latest Handback DateTime
    $ws.Range("K6").Value = $handbackDateTime

    # Error Detail
    $ws.Range("P6").Value = $errorDetail

    $ws.Hyperlinks.Add($ws.Range("I6"), $targetUrl, "", "", "af758db3-7343-497e-8cc8-91cff17baa62.md") | Out-Null
}

Update-LocaleSheet "zh-cn" "2016-10-14 07:33:58"
Update-LocaleSheet "de-de" "2016-10-14 07:34:15"
